$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns L and M
$ws.Range("L1").Value = "prf_chg_attn_fac"
$ws.Range("M1").Value = "pu_prf_pass0"

# New column width for column L (raw OOXML width of 15 once Excel's
# character-width padding of 5/7 is applied by ColumnWidth)
$ws.Columns("L").ColumnWidth = 14.285714285714286

# New data values for columns L and M, rows 2-8
$lValues = @(1.1599999999999999, 2.6, 11.34, 20.57, 29.3, 30.29, 30.23)
$mValues = @(0.0063, 0.0063, 0.0063, 0.0063, 0.0063, 0.0063, 0.0063)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
    $ws.Cells.Item($row, 13).Value = $mValues[$i]
}

# Update the selected cell to match the final workbook state
$ws.Range("J16").Select()
